# Refresh computed market-profit columns (H:N) across the Leve-profit sheets.
# Mirrors a scheduled data-refresh run: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ)
# and LeveProfit(NQ/HQ) are recomputed per row from updated market-board snapshots.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 275
$ws.Range("I41").Value = 306.625
$ws.Range("J41").Value = 243.375
$ws.Range("K41").Value = 306.625
$ws.Range("L41").Value = 243.375
$ws.Range("M41").Value = 133.375
$ws.Range("N41").Value = -1123.375
$ws.Range("H62").Value = 2549.5833
$ws.Range("I62").Value = 2508.6365
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2508.6365
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1884.6365
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2549.5833
$ws.Range("I65").Value = 2508.6365
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 12543.1825
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -9423.182500000001
$ws.Range("N65").Value = -21240
$ws.Range("H99").Value = 579
$ws.Range("I99").Value = 579
$ws.Range("K99").Value = 1737
$ws.Range("M99").Value = -239
$ws.Range("H129").Value = 764.8148
$ws.Range("J129").Value = 1044.5834
$ws.Range("L129").Value = 3133.7502
$ws.Range("N129").Value = -13133.7502
$ws.Range("H132").Value = 559045.5600000001
$ws.Range("I132").Value = 1967.6624
$ws.Range("K132").Value = 5902.9872
$ws.Range("M132").Value = -3372.9872
$ws.Range("H135").Value = 30675.371
$ws.Range("I135").Value = 37310.785
$ws.Range("J135").Value = 4133.7144
$ws.Range("K135").Value = 335797.0650000001
$ws.Range("L135").Value = 37203.4296
$ws.Range("M135").Value = -333262.0650000001
$ws.Range("N135").Value = -42273.4296
$ws.Range("H137").Value = 1889651
$ws.Range("I137").Value = 2327686
$ws.Range("J137").Value = 6100.7
$ws.Range("K137").Value = 6983058
$ws.Range("L137").Value = 18302.1
$ws.Range("M137").Value = -6980508
$ws.Range("N137").Value = -23402.1
$ws.Range("H138").Value = 4002663
$ws.Range("I138").Value = 2212.3462
$ws.Range("J138").Value = 8336484.5
$ws.Range("K138").Value = 6637.0386
$ws.Range("L138").Value = 25009453.5
$ws.Range("M138").Value = -1497.0386
$ws.Range("N138").Value = -25019733.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3654.31
$ws.Range("I32").Value = 3626.8572
$ws.Range("J32").Value = 4999.5
$ws.Range("K32").Value = 3626.8572
$ws.Range("L32").Value = 4999.5
$ws.Range("M32").Value = -3339.8572
$ws.Range("N32").Value = -5573.5
$ws.Range("H61").Value = 111335656
$ws.Range("I61").Value = 143002620
$ws.Range("J61").Value = 501250
$ws.Range("K61").Value = 143002620
$ws.Range("L61").Value = 501250
$ws.Range("M61").Value = -143002408
$ws.Range("N61").Value = -501674
$ws.Range("H74").Value = 5083198.5
$ws.Range("I74").Value = 6557146.5
$ws.Range("J74").Value = 71774.266
$ws.Range("K74").Value = 6557146.5
$ws.Range("L74").Value = 71774.266
$ws.Range("M74").Value = -6556272.5
$ws.Range("N74").Value = -73522.266
$ws.Range("H77").Value = 5083198.5
$ws.Range("I77").Value = 6557146.5
$ws.Range("J77").Value = 71774.266
$ws.Range("K77").Value = 32785732.5
$ws.Range("L77").Value = 358871.33
$ws.Range("M77").Value = -32781364.5
$ws.Range("N77").Value = -367607.33
$ws.Range("H132").Value = 9807671
$ws.Range("I132").Value = 11629541
$ws.Range("J132").Value = 15119.75
$ws.Range("K132").Value = 34888623
$ws.Range("L132").Value = 45359.25
$ws.Range("M132").Value = -34886093
$ws.Range("N132").Value = -50419.25
$ws.Range("H136").Value = 111335656
$ws.Range("I136").Value = 143002620
$ws.Range("J136").Value = 501250
$ws.Range("K136").Value = 429007860
$ws.Range("L136").Value = 1503750
$ws.Range("M136").Value = -429005310
$ws.Range("N136").Value = -1508850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8134498.5
$ws.Range("I134").Value = 4811.484
$ws.Range("K134").Value = 14434.452
$ws.Range("M134").Value = -11899.452

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 42979.816
$ws.Range("I31").Value = 34194.258
$ws.Range("J31").Value = 54821.22
$ws.Range("K31").Value = 34194.258
$ws.Range("L31").Value = 54821.22
$ws.Range("M31").Value = -33899.258
$ws.Range("N31").Value = -55411.22
$ws.Range("H34").Value = 42979.816
$ws.Range("I34").Value = 34194.258
$ws.Range("J34").Value = 54821.22
$ws.Range("K34").Value = 34194.258
$ws.Range("L34").Value = 54821.22
$ws.Range("M34").Value = -33992.258
$ws.Range("N34").Value = -55225.22
$ws.Range("H58").Value = 15626607
$ws.Range("I58").Value = 28572784
$ws.Range("K58").Value = 28572784
$ws.Range("M58").Value = -28572581
$ws.Range("H94").Value = 3723.9546
$ws.Range("I94").Value = 10400.333
$ws.Range("J94").Value = 1220.3125
$ws.Range("K94").Value = 10400.333
$ws.Range("L94").Value = 1220.3125
$ws.Range("M94").Value = -9949.333000000001
$ws.Range("N94").Value = -2122.3125
$ws.Range("H103").Value = 27666.666
$ws.Range("I103").Value = 27666.666
$ws.Range("K103").Value = 27666.666
$ws.Range("M103").Value = -26494.666
$ws.Range("H122").Value = 2296.762
$ws.Range("I122").Value = 1700.2667
$ws.Range("J122").Value = 3788
$ws.Range("K122").Value = 5100.800099999999
$ws.Range("L122").Value = 11364
$ws.Range("M122").Value = -2650.800099999999
$ws.Range("N122").Value = -16264
$ws.Range("H132").Value = 27552.564
$ws.Range("I132").Value = 1615.1034
$ws.Range("J132").Value = 102771.2
$ws.Range("K132").Value = 4845.3102
$ws.Range("L132").Value = 308313.6
$ws.Range("M132").Value = -2315.3102
$ws.Range("N132").Value = -313373.6
$ws.Range("H136").Value = 15626607
$ws.Range("I136").Value = 28572784
$ws.Range("K136").Value = 85718352
$ws.Range("M136").Value = -85715802

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4183.706
$ws.Range("I64").Value = 2108
$ws.Range("J64").Value = 4628.5
$ws.Range("K64").Value = 6324
$ws.Range("L64").Value = 13885.5
$ws.Range("M64").Value = -6054
$ws.Range("N64").Value = -14425.5
$ws.Range("H67").Value = 4183.706
$ws.Range("I67").Value = 2108
$ws.Range("J67").Value = 4628.5
$ws.Range("K67").Value = 6324
$ws.Range("L67").Value = 13885.5
$ws.Range("M67").Value = -5388
$ws.Range("N67").Value = -15757.5
$ws.Range("H114").Value = 10753321
$ws.Range("I114").Value = 498.0909
$ws.Range("J114").Value = 16667373
$ws.Range("K114").Value = 1494.2727
$ws.Range("L114").Value = 50002119
$ws.Range("M114").Value = 1759.7273
$ws.Range("N114").Value = -50008627
$ws.Range("H117").Value = 13333899
$ws.Range("J117").Value = 66666664
$ws.Range("L117").Value = 199999992
$ws.Range("N117").Value = -200006876
$ws.Range("H121").Value = 65861330
$ws.Range("I121").Value = 377.5
$ws.Range("J121").Value = 92205704
$ws.Range("K121").Value = 1132.5
$ws.Range("L121").Value = 276617112
$ws.Range("M121").Value = 177.5
$ws.Range("N121").Value = -276619732
$ws.Range("H131").Value = 821.5
$ws.Range("J131").Value = 919.1429000000001
$ws.Range("L131").Value = 2757.4287
$ws.Range("N131").Value = -12837.4287
$ws.Range("H136").Value = 3148.739
$ws.Range("I136").Value = 2780.7273
$ws.Range("J136").Value = 3486.0833
$ws.Range("K136").Value = 8342.1819
$ws.Range("L136").Value = 10458.2499
$ws.Range("M136").Value = -3242.1819
$ws.Range("N136").Value = -20658.2499
$ws.Range("H137").Value = 2018.5
$ws.Range("I137").Value = 963.75
$ws.Range("J137").Value = 2229.45
$ws.Range("K137").Value = 2891.25
$ws.Range("L137").Value = 6688.349999999999
$ws.Range("M137").Value = 2208.75
$ws.Range("N137").Value = -16888.35
$ws.Range("H141").Value = 4064.8
$ws.Range("I141").Value = 4064.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 12194.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -7014.400000000001
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2030
$ws.Range("I126").Value = 1465.2
$ws.Range("J126").Value = 2971.3333
$ws.Range("K126").Value = 4395.6
$ws.Range("L126").Value = 8913.999899999999
$ws.Range("M126").Value = -1925.6
$ws.Range("N126").Value = -13853.9999
$ws.Range("H132").Value = 52129.2
$ws.Range("I132").Value = 36147.1
$ws.Range("J132").Value = 94263.82000000001
$ws.Range("K132").Value = 108441.3
$ws.Range("L132").Value = 282791.46
$ws.Range("M132").Value = -105911.3
$ws.Range("N132").Value = -287851.46

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20867.303
$ws.Range("I132").Value = 2275.8975
$ws.Range("J132").Value = 63518.176
$ws.Range("K132").Value = 6827.6925
$ws.Range("L132").Value = 190554.528
$ws.Range("M132").Value = -4297.6925
$ws.Range("N132").Value = -195614.528
$ws.Range("H136").Value = 32811.85
$ws.Range("I136").Value = 18957.967
$ws.Range("J136").Value = 204600
$ws.Range("K136").Value = 56873.901
$ws.Range("L136").Value = 613800
$ws.Range("M136").Value = -54323.901
$ws.Range("N136").Value = -618900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H132").Value = 37771.656
$ws.Range("I132").Value = 27947.107
$ws.Range("J132").Value = 57966.555
$ws.Range("K132").Value = 83841.321
$ws.Range("L132").Value = 173899.665
$ws.Range("M132").Value = -81311.321
$ws.Range("N132").Value = -178959.665
$ws.Range("H136").Value = 39301.758
$ws.Range("I136").Value = 24439.07
$ws.Range("J136").Value = 97401.37
$ws.Range("K136").Value = 73317.20999999999
$ws.Range("L136").Value = 292204.11
$ws.Range("M136").Value = -70767.20999999999
$ws.Range("N136").Value = -297304.11
